$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mobile Number (C2): keep as text with a leading zero preserved, the same
# way Excel stores a manually apostrophe-prefixed entry (quotePrefix style)
# rather than converting it back to a number.
$ws.Range("C2").Value = "'0508037025"

# Assured Name (B2): updated client name.
$ws.Range("B2").Value = "SALEM KHALIFA SALEM RASHED"

# Match the saved selection/active cell shown in the workbook.
[void]$ws.Range("B2").Select()
